# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value2  = 3056
$wsExhibit.Range("F3").Value2  = 474
$wsExhibit.Range("F5").Value2  = 39
$wsExhibit.Range("F7").Value2  = 1038
$wsExhibit.Range("F8").Value2  = 14714
$wsExhibit.Range("F9").Value2  = 173
$wsExhibit.Range("F11").Value2 = 5864
$wsExhibit.Range("F12").Value2 = 598
$wsExhibit.Range("F15").Value2 = 77
$wsExhibit.Range("F16").Value2 = 1240
$wsExhibit.Range("F18").Value2 = 90
$wsExhibit.Range("F21").Value2 = 2945
$wsExhibit.Range("F22").Value2 = 90
$wsExhibit.Range("F23").Value2 = 10666
$wsExhibit.Range("F26").Value2 = 99
$wsExhibit.Range("F27").Value2 = 3745
$wsExhibit.Range("F28").Value2 = 249

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value2  = 3056
$wsAll.Range("F4").Value2  = 474
$wsAll.Range("F6").Value2  = 39
$wsAll.Range("F8").Value2  = 1038
$wsAll.Range("F9").Value2  = 14714
$wsAll.Range("F10").Value2 = 173
$wsAll.Range("F12").Value2 = 5864
$wsAll.Range("F13").Value2 = 598
$wsAll.Range("F16").Value2 = 77
$wsAll.Range("F17").Value2 = 1240
$wsAll.Range("F19").Value2 = 90
$wsAll.Range("F22").Value2 = 2945
$wsAll.Range("F23").Value2 = 90
$wsAll.Range("F25").Value2 = 10666
$wsAll.Range("F28").Value2 = 99
$wsAll.Range("F29").Value2 = 3745
$wsAll.Range("F30").Value2 = 249
